$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The DICIEMBRE (December) row in each region block was left at 0/placeholder
# values. Bring each one into line with the NOVIEMBRE (November) row directly
# above it - same fixed_charge / cost_verano / cost_no_verano figures and the
# same cell formatting (including the special highlighted style used for the
# first region's row).

$targets = @(
    @{ Src = 12; Dst = 13; D = 98.04000000000001; E = 4.201; F = 3.61 },
    @{ Src = 24; Dst = 25; D = 98.04000000000001; E = 4.58;  F = 3.61 },
    @{ Src = 36; Dst = 37; D = 98.04000000000001; E = 4.32;  F = 4.32 },
    @{ Src = 48; Dst = 49; D = 98.04000000000001; E = 4.05;  F = 4.05 },
    @{ Src = 60; Dst = 61; D = 98.04000000000001; E = 3.95;  F = 3.95 },
    @{ Src = 72; Dst = 73; D = 98.04000000000001; E = 4.01;  F = 4.01 }
)

foreach ($t in $targets) {
    $srcRange = "D" + $t.Src + ":F" + $t.Src
    $dstRange = "D" + $t.Dst + ":F" + $t.Dst

    # Copy the formatting (styles) from the source row onto the target row.
    $ws.Range($srcRange).Copy()
    $ws.Range($dstRange).PasteSpecial(-4122)

    # Now write the actual values onto the target row.
    $ws.Range("D" + $t.Dst).Value = $t.D
    $ws.Range("E" + $t.Dst).Value = $t.E
    $ws.Range("F" + $t.Dst).Value = $t.F
}
